# Auto-generated PowerShell COM-interop script
# Updates the cryptos list (Coin/Link/Price/Volume) to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.849.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.898.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "

$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5017"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2970"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06817"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.87%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.907.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07321"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "91.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.091"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6771"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.824.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.152.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.866"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "181.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +33.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.074"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.341"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.941"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.395"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.19%  "

$ws.Range("E30").Value = "  +4.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08984"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05262"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7447"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.40%  "

$ws.Range("E35").Value = "  +2.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.670"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01934"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +17.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.720"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.178"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9369"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4384"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.823"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.24%  "

$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.730"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1344"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05842"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.25%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.3914"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.39%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.532"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.01%  "

$ws.Range("E50").Value = "  +2.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.382"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.53%  "
